$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fancycopy")

# Bug fix 1: rename the generic placeholder hyperlink text "t" (row 7) to the
# clearer label "link".
$ws.Range("C7").Value2 = "link"

# Bug fix 2: the second hyperlink row (row 8) was incorrectly reusing the
# same display text as row 7 ("t") even though it points to a different
# link; give it its own distinct label "link2" instead.
$ws.Range("C8").Value2 = "link2"
